$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.050.11"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "3.066.32"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("E9").Value = "  +2.76%  "
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").Value = "3.591.65"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "58.085.05"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").Value = "3.058.47"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("E18").Value = "  +3.28%  "
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.55%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "0.0₃0906"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  +7.12%  "
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("E31").Value = "  +4.18%  "
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("E37").Value = "  +3.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0681"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.92%  "
$ws.Range("D39").Value = "3.110.10"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.659"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "2.294.86"
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("E45").Value = "  +5.99%  "
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.77%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.944"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.730"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0879"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.65%  "
